# Add 2022-Q4 data
#
# The workbook currently has sheets: 总计, 2022-Q3, 2022-Q1
# After the edit it should have: 总计, 2022-Q4, 2022-Q3, 2022-Q1
# - A new "2022-Q4" sheet is inserted (as a copy of "2022-Q3", with updated figures)
# - The "总计" (summary) sheet gets a new row for 2022-Q4, pushing the existing
#   2022-Q3 / 2022-Q1 rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Build the new "2022-Q4" worksheet by duplicating "2022-Q3" (this keeps all
#    styles / column layout identical) and placing it immediately before the
#    existing "2022-Q3" sheet.
# ---------------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Copy($wsQ3, $null)
$wsQ4 = $wb.Worksheets.Item("2022-Q3 (2)")
$wsQ4.Name = "2022-Q4"

# Update the figures on the new 2022-Q4 sheet with the latest quarter's values.
# Columns D:G hold numbers formatted as text in this workbook, so force a text
# number format before assigning them so they keep their original data type.
$wsQ4.Range("D2:G3").NumberFormat = "@"

$wsQ4.Range("D2").Value = "2.33"
$wsQ4.Range("E2").Value = "93.03"
$wsQ4.Range("F2").Value = "4.81"
$wsQ4.Range("G2").Value = "0.1121"
$wsQ4.Range("H2").Value = 8

$wsQ4.Range("D3").Value = "0.29"
$wsQ4.Range("E3").Value = "93.03"
$wsQ4.Range("F3").Value = "4.81"
$wsQ4.Range("G3").Value = "0.0139"
$wsQ4.Range("H3").Value = 8

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert the new 2022-Q4 row at the top
#    of the data (row 2), pushing the old 2022-Q3 / 2022-Q1 rows down one row.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Row 4 (2022-Q1) is brand new territory for this sheet, so first copy the
# formatting of an existing data row (row 2) onto it before filling in values.
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A4").PasteSpecial(-4122)

# Fill bottom-up so we never clobber a row before reading its old value.
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2022-Q1"
$wsTotal.Range("C4").Value = 3
$wsTotal.Range("D4").Value = 0.12

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.12

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.13

# ---------------------------------------------------------------------------
# 3. Restore the originally selected tab ("2022-Q1", last sheet) so the
#    active-sheet bookkeeping matches the source workbook.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q1").Activate()
